$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.177.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.38%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.310.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.74%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.27%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.60%  "

# Row 6
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.67%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.49%  "

# Row 8
$ws.Range("E8").Value = "  -0.32%  "

# Row 9
$ws.Range("E9").Value = "  -0.77%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.93%  "

# Row 11
$ws.Range("E11").Value = "  -0.51%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.40%  "

# Row 13
$ws.Range("E13").Value = "  +17.23%  "

# Row 14
$ws.Range("E14").Value = "  -0.55%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.33%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.656.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.69%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.300.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.09%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.155.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.64%  "

# Row 19
$ws.Range("E19").Value = "  -0.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.73%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.11%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.90%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "255.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.84%  "

# Row 25
$ws.Range("E25").Value = "  -5.39%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.09%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.56%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.96%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.09%  "

# Row 32
$ws.Range("E32").Value = "  -2.83%  "

# Row 33
$ws.Range("E33").Value = "  -0.49%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.81%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.76%  "

# Row 36
$ws.Range("E36").Value = "  +1.13%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.88%  "

# Row 38
$ws.Range("E38").Value = "  +0.88%  "

# Row 39
$ws.Range("E39").Value = "  -1.46%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.86%  "

# Row 42
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.67%  "

# Row 43
$ws.Range("E43").Value = "  -3.39%  "

# Row 44
$ws.Range("E44").Value = "  -0.04%  "

# Row 45
$ws.Range("E45").Value = "  -7.46%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.89%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "108.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.57%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.84%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.12%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0990"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.14%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.12%  "
